# engineer_data.py change: capitalize the hex byte codes in the
# "doip" (G) and "uds" (H) columns so they are easier to look up,
# e.g. "0x02:0xfd:0x00" -> "0x02:0xFD:0x00". The leading "0x" marker
# itself stays lowercase; only the two hex digits after it are
# uppercased. Cells that are not colon-separated hex byte strings
# (e.g. "N/A") are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function ConvertTo-UpperHexCode($s) {
    $segments = $s.Split(":")
    $result = ""
    foreach ($seg in $segments) {
        if ($result -ne "") {
            $result = $result + ":"
        }
        if ($seg.Length -ge 2 -and $seg.Substring(0, 2).Equals("0x")) {
            $result = $result + "0x" + $seg.Substring(2).ToUpper()
        } else {
            $result = $result + $seg
        }
    }
    return $result
}

$lastRow = 33
for ($row = 2; $row -le $lastRow; $row++) {
    foreach ($col in @("G", "H")) {
        $cell = $ws.Range($col + $row)
        $txt = $cell.Text
        if ($txt -ne $null -and $txt.Contains("0x")) {
            $cell.Value = ConvertTo-UpperHexCode($txt)
        }
    }
}
